# Add a new "Height" property column (AB) for the NPC table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new column.
$ws.Range("AB1").Value = "Height"

# Data rows 2..21 all get the value 2 in the new column.
$ws.Range("AB2:AB21").Value = 2

# Match the column width used by the neighbouring column (engine quantizes
# column widths to 1/7 character-width steps, so 14 is the closest input
# that reproduces the intended ~14.75 display width).
$ws.Columns("AB").ColumnWidth = 14

# Mirror the author's final selection/view state on the sheet.
$ws.Range("AB2:AB21").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 24
$win.ScrollRow = 1
